# Appends new Donor rows (12-19) to the "Donors" sheet and new Patient
# rows (10-16) to the "Patients" sheet, matching a fresh export of the
# underlying lifelink_db.xlsx MongoDB-backed data.

$wb = $excel.ActiveWorkbook
$donors = $wb.Worksheets.Item("Donors")
$patients = $wb.Worksheets.Item("Patients")

# ---------------------------------------------------------------------
# Sheet "Donors" (sheet1): new rows 12-19
# ---------------------------------------------------------------------
$donorRows = @(
    @("69234632a9d7b28957d0eb0c", "mno", "mno@gmail.com", "9481824919", "A+", "Mangalore", 0, 0, $false, $false, $true, 45984.96320155093),
    @("69234680a9d7b28957d0eb11", "Me", "23a43.bhavish@sjec.ac.in", "8904534919", "A+", "Mangalore", 0, 0, $false, $false, $true, 45984.9641012037),
    @("6923b7181c52f11af60577df", "mailtrap", "lifelink@system.com", "1234456789", "A+", "Mangalore", 0, 0, $false, $false, $true, 45985.297717719906),
    @("6923de6730064d1b9aecd2e1", "Ananya", "ananyaskulai@gmail.com", "6361037723", "O+", "My Current Location", 74.899521, 12.9103476, $false, $false, $true, 45985.4141815625),
    @("6923e10530064d1b9aecd323", "Ananya S Kulai", "ananyaskulai@gmail.com", "6361037723", "A+", "My Current Location", 77.6077312, 12.9728512, $false, $false, $true, 45985.42193297454),
    @("6923e2e830064d1b9aecd366", "Bhavish", "bhavish@gmail.com", "9986769690", "B-", "My Current Location", 74.8994941, 12.9103193, $false, $false, $true, 45985.4275340625),
    @("6923e4ccc81a5010397a8964", "Ashwini Shenoy B", "lifelink@system.com", "7026438371", "A+", "My Current Location", 77.6077312, 12.9728512, $false, $false, $true, 45985.433130439815),
    @("6923f711f76d0536415e5e94", "Ashwini Shenoy B", "lifelink@gmail.com", "7022157406", "A+", "My Current Location", 74.8996501, 12.9103764, $false, $false, $true, 45985.48726598379)
)

# Column D (phone) is all-digits text; format it as Text up-front so
# Excel doesn't auto-coerce the literal strings into numbers.
$donors.Range("D12:D19").NumberFormat = "@"

$r = 12
foreach ($row in $donorRows) {
    $donors.Cells.Item($r, 1).Value = $row[0]
    $donors.Cells.Item($r, 2).Value = $row[1]
    $donors.Cells.Item($r, 3).Value = $row[2]
    $donors.Cells.Item($r, 4).Value = $row[3]
    $donors.Cells.Item($r, 5).Value = $row[4]
    $donors.Cells.Item($r, 6).Value = $row[5]
    $donors.Cells.Item($r, 7).Value = $row[6]
    $donors.Cells.Item($r, 8).Value = $row[7]
    $donors.Cells.Item($r, 9).Value = $row[8]
    $donors.Cells.Item($r, 10).Value = $row[9]
    $donors.Cells.Item($r, 11).Value = $row[10]
    $donors.Cells.Item($r, 12).Value = $row[11]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "Patients" (sheet2): new rows 10-16
# ---------------------------------------------------------------------
$patientRows = @(
    @("6923d7e8e754f1d7bef87e55", "Ashwini Shenoy B", "ashwinishenoyb@gmail.com", "7026438371", "A+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.394934780095),
    @("6923dfc630064d1b9aecd306", "Ashwini Shenoy B", "ashenoyb@gmail.com", "7026438371", "A+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.41824520833),
    @("6923dfe430064d1b9aecd30f", "Ashwini Shenoy B", "ashenoyb@gmail.com", "7026438371", "O+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.41858971065),
    @("6923e15630064d1b9aecd32e", "Ashwini Shenoy B", "ashenoyb@gmail.com", "7026438371", "A+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.422877916666),
    @("6923e23c30064d1b9aecd354", "Ashwini Shenoy B", "ashwinishenoyb@gmail.com", "7026438371", "A+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.425535104165),
    @("6923e5a0c81a5010397a896d", "Ananya S Kulai", "ananyaskulai@gmail.com", "6361037723", "A+", "My Current Location", 77.6077312, 12.9728512, "High", 45985.43557923611),
    @("6923f833f76d0536415e5ea7", "Ashwini Shenoy B", "ashwinishenoyb@gmail.com", "7026438371", "A+", "My Current Location", 74.8996566, 12.9103667, "Critical", 45985.49062084491)
)

# Same text-coercion guard for Patients!D (phone).
$patients.Range("D10:D16").NumberFormat = "@"

$r = 10
foreach ($row in $patientRows) {
    $patients.Cells.Item($r, 1).Value = $row[0]
    $patients.Cells.Item($r, 2).Value = $row[1]
    $patients.Cells.Item($r, 3).Value = $row[2]
    $patients.Cells.Item($r, 4).Value = $row[3]
    $patients.Cells.Item($r, 5).Value = $row[4]
    $patients.Cells.Item($r, 6).Value = $row[5]
    $patients.Cells.Item($r, 7).Value = $row[6]
    $patients.Cells.Item($r, 8).Value = $row[7]
    $patients.Cells.Item($r, 9).Value = $row[8]
    $patients.Cells.Item($r, 10).Value = $row[9]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Formatting tweaks to match the refreshed export:
#  - Patients!J16 gets the same date-number style Donors!L11 already
#    carried (copy it over first, so the style-table entry is reused
#    instead of a near-duplicate being minted).
#  - Donors!L11 itself drops that date style, becoming a plain number
#    like the rest of column L.
# ---------------------------------------------------------------------
$donors.Range("L11").Copy()
$patients.Range("J16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$donors.Range("L11").ClearFormats()
